# Updates cryptocurrency price/volume figures in the active worksheet
# to reflect the latest scrape, per the automated GitHub Actions job.
# Price values (column D) are stored as plain text (they use a
# locale-style grouping like "27.331.69"), so force a text number
# format before assigning to stop Excel from reinterpreting them as
# numeric/date values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.331.69"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.711.21"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.20"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5294"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06663"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.84"
$ws.Range("E10").Value = "  -3.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07694"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.497"
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.946.98"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.718.79"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5826"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8208"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.00"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.356.38"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "221.20"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.626"
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.44"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.003"
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.74"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.689"
$ws.Range("E26").Value = "  -2.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1204"
$ws.Range("E27").Value = "  -2.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.228"
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.24"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05332"
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.461"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.433"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.634"
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.879"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9510"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.394"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5832"
$ws.Range("E38").Value = "  -2.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01633"
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.114.09"
$ws.Range("E40").Value = "  +5.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.804"
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8389"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.20"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.854.00"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.62"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("E48").Value = "  +2.22%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.080"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05226"
$ws.Range("E51").Value = "  -0.34%  "
